$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=28834; B="Sr. João Vitor Camargo";  C="P&D";               D="Viagem de negocios"; E=3; F=45088; G=9270.57}
    @{Row=3;  A=94421; B="Aurora Sales";             C="Marketing";         D="Doenca";             E=2; F=45084; G=7347.95}
    @{Row=4;  A=21118; B="Dr. Ravi Lucca Mendes";    C="TI";                D="Doenca";             E=8; F=45084; G=4523.79}
    @{Row=5;  A=8913;  B="Ana Liz Macedo";           C="Engenharia";        D="Outros";             E=7; F=45078; G=4555.55}
    @{Row=6;  A=21337; B="Dr. Luiz Otávio Sá";       C="Recursos Humanos";  D="Outros";             E=7; F=45088; G=8911.299999999999}
    @{Row=7;  A=32264; B="Sr. Samuel Sales";         C="Vendas";            D="Doenca";             E=3; F=45100; G=2935.28}
    @{Row=8;  A=6229;  B="Laura das Neves";          C="Vendas";            D="Problemas pessoais"; E=2; F=45098; G=8186.67}
    @{Row=9;  A=69262; B="Erick Garcia";             C="Operacoes";         D="Problemas pessoais"; E=5; F=45089; G=8186.37}
    @{Row=10; A=75690; B="Dr. Noah Barros";          C="Vendas";            D="Viagem de negocios"; E=7; F=45095; G=6039.31}
    @{Row=11; A=74867; B="Bernardo Nascimento";      C="TI";                D="Doenca";             E=7; F=45100; G=6813.51}
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = $rec.A
    $ws.Cells.Item($r, 2).Value = $rec.B
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
}
